$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = "'27.013.53"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +0.99%  "
$ws.Range('E2').Style = 'Normal'
# Row 3
$ws.Range('D3').Value = "'1.644.98"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +0.05%  "
$ws.Range('E3').Style = 'Normal'
# Row 4
$ws.Range('E4').Value = "'  -0.32%  "
$ws.Range('E4').Style = 'Normal'
# Row 5
$ws.Range('D5').Value = "'219.76"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +1.60%  "
$ws.Range('E5').Style = 'Normal'
# Row 6
$ws.Range('D6').Value = "'0.499"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -0.33%  "
$ws.Range('E6').Style = 'Normal'
# Row 7
$ws.Range('E7').Value = "'  -0.27%  "
$ws.Range('E7').Style = 'Normal'
# Row 8
$ws.Range('D8').Value = "'0.253"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  +0.34%  "
$ws.Range('E8').Style = 'Normal'
# Row 9
$ws.Range('E9').Value = "'  -0.39%  "
$ws.Range('E9').Style = 'Normal'
# Row 10
$ws.Range('D10').Value = "'19.44"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  +1.42%  "
$ws.Range('E10').Style = 'Normal'
# Row 11
$ws.Range('D11').Value = "'0.0848"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  +0.61%  "
$ws.Range('E11').Style = 'Normal'
# Row 12
$ws.Range('D12').Value = "'1.874.17"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  +0.03%  "
$ws.Range('E12').Style = 'Normal'
# Row 13
$ws.Range('D13').Value = "'1.634.28"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -0.56%  "
$ws.Range('E13').Style = 'Normal'
# Row 14
$ws.Range('D14').Value = "'4.18"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +0.24%  "
$ws.Range('E14').Style = 'Normal'
# Row 15
$ws.Range('E15').Value = "'  +0.57%  "
$ws.Range('E15').Style = 'Normal'
# Row 16
$ws.Range('D16').Value = "'66.02"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +2.09%  "
$ws.Range('E16').Style = 'Normal'
# Row 17
$ws.Range('D17').Value = "'26.986.59"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +0.89%  "
$ws.Range('E17').Style = 'Normal'
# Row 18
$ws.Range('E18').Value = "'  -0.15%  "
$ws.Range('E18').Style = 'Normal'
# Row 19
$ws.Range('D19').Value = "'218.49"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  +2.04%  "
$ws.Range('E19').Style = 'Normal'
# Row 20
$ws.Range('E20').Value = "'  -0.24%  "
$ws.Range('E20').Style = 'Normal'
# Row 21
$ws.Range('D21').Value = "'4.41"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  -0.14%  "
$ws.Range('E21').Style = 'Normal'
# Row 22
$ws.Range('E22').Value = "'  +6.15%  "
$ws.Range('E22').Style = 'Normal'
# Row 23
$ws.Range('D23').Value = "'2.44"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  -1.56%  "
$ws.Range('E23').Style = 'Normal'
# Row 24
$ws.Range('E24').Value = "'  -1.23%  "
$ws.Range('E24').Style = 'Normal'
# Row 25
$ws.Range('D25').Value = "'148.30"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  +1.59%  "
$ws.Range('E25').Style = 'Normal'
# Row 26
$ws.Range('E26').Value = "'  -0.30%  "
$ws.Range('E26').Style = 'Normal'
# Row 27
$ws.Range('D27').Value = "'7.32"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  +2.56%  "
$ws.Range('E27').Style = 'Normal'
# Row 28
$ws.Range('D28').Value = "'0.119"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  -0.02%  "
$ws.Range('E28').Style = 'Normal'
# Row 29
$ws.Range('D29').Value = "'15.83"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  +1.26%  "
$ws.Range('E29').Style = 'Normal'
# Row 30
$ws.Range('E30').Value = "'  +0.52%  "
$ws.Range('E30').Style = 'Normal'
# Row 31
$ws.Range('E31').Value = "'  +1.40%  "
$ws.Range('E31').Style = 'Normal'
# Row 32
$ws.Range('E32').Value = "'  +0.36%  "
$ws.Range('E32').Style = 'Normal'
# Row 33
$ws.Range('D33').Value = "'3.00"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  -0.12%  "
$ws.Range('E33').Style = 'Normal'
# Row 34
$ws.Range('E34').Value = "'  +1.68%  "
$ws.Range('E34').Style = 'Normal'
# Row 35
$ws.Range('D35').Value = "'1.270.84"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  -1.23%  "
$ws.Range('E35').Style = 'Normal'
# Row 36
$ws.Range('E36').Value = "'  -0.04%  "
$ws.Range('E36').Style = 'Normal'
# Row 37
$ws.Range('E37').Value = "'  -1.90%  "
$ws.Range('E37').Style = 'Normal'
# Row 38
$ws.Range('D38').Value = "'0.533"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  -0.61%  "
$ws.Range('E38').Style = 'Normal'
# Row 39
$ws.Range('E39').Value = "'  +0.69%  "
$ws.Range('E39').Style = 'Normal'
# Row 40
$ws.Range('E40').Value = "'  -0.22%  "
$ws.Range('E40').Style = 'Normal'
# Row 41
$ws.Range('D41').Value = "'0.807"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  +0.36%  "
$ws.Range('E41').Style = 'Normal'
# Row 42
$ws.Range('D42').Value = "'5.36"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  +0.94%  "
$ws.Range('E42').Style = 'Normal'
# Row 43
$ws.Range('D43').Value = "'1.784.53"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  +0.06%  "
$ws.Range('E43').Style = 'Normal'
# Row 44
$ws.Range('D44').Value = "'61.87"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  +0.01%  "
$ws.Range('E44').Style = 'Normal'
# Row 45
$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D45').Value = "'92.53"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  +0.90%  "
$ws.Range('E45').Style = 'Normal'
# Row 46
$ws.Range('B46').Value = 'MXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D46').Value = "'2.05"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -7.95%  "
$ws.Range('E46').Style = 'Normal'
# Row 47
$ws.Range('E47').Value = "'  -0.31%  "
$ws.Range('E47').Style = 'Normal'
# Row 48
$ws.Range('E48').Value = "'  -0.97%  "
$ws.Range('E48').Style = 'Normal'
# Row 49
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = "'7.60"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  -1.18%  "
$ws.Range('E49').Style = 'Normal'
# Row 50
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').Value = "'0.0970"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  -0.14%  "
$ws.Range('E50').Style = 'Normal'
# Row 51
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').Value = "'0.406"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +0.04%  "
$ws.Range('E51').Style = 'Normal'
